# Edit script: "Menu has background. Some more animations have been added/fixed."
#
# Summary of changes applied:
#  - Row 20, col A: "main_menu.png" -> "bg3.png"
#  - Rows 25 and 26: all cell contents cleared (values removed, styles kept)
#  - Rows 29 and 30, col D: cell text changed from the raw URL to "Original"
#    (hyperlink target kept the same; Excel will record the old URL text as
#    the hyperlink's display attribute automatically)
#  - Rows 29 and 30, col E: "Finished" -> "Placeholder"
#  - Selection / scroll position of the sheet view moved to A17

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: main menu picture filename now points at the new background file
$ws.Range("A20").Value = "bg3.png"

# --- Rows 25 & 26: clear all 5 columns (A:E) leaving the cells blank
$ws.Range("A25:E25").ClearContents()
$ws.Range("A26:E26").ClearContents()

# --- Rows 29 & 30: the hyperlink text is replaced with "Original" while the
#     hyperlink itself (and its target) stays intact, and the status moves
#     from Finished to Placeholder
$ws.Range("D29").Value = "Original"
$ws.Range("E29").Value = "Placeholder"

$ws.Range("D30").Value = "Original"
$ws.Range("E30").Value = "Placeholder"

# --- Move the view / active selection to A17
$ws.Range("A17").Select()
